# Update the "K" column (column G) values for rows 2-11 on Sheet1,
# per the regenerated save_data (K instead of Strike#, std/mean, s_vals).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 3
    6  = 0
    7  = 1
    8  = 1
    9  = 3
    10 = 1
    11 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
